# Update Esparragos weekly price rows (12-26) and append two new rows (27-28)
# for "Hortaliza, Vega Modelo de Temuco - Esparragos".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Cells.Item(12,4).Value = 44483
$ws.Cells.Item(12,9).Value = "Extra"
$ws.Cells.Item(12,10).Value = 50
$ws.Cells.Item(12,11).Value = 2000
$ws.Cells.Item(12,12).Value = 2000
$ws.Cells.Item(12,13).Value = 2000
$ws.Cells.Item(12,16).Value = 2000

# Row 13
$ws.Cells.Item(13,4).Value = 44483
$ws.Cells.Item(13,10).Value = 500
$ws.Cells.Item(13,11).Value = 1300
$ws.Cells.Item(13,12).Value = 1500
$ws.Cells.Item(13,13).Value = 1420
$ws.Cells.Item(13,16).Value = 1420

# Row 14
$ws.Cells.Item(14,4).Value = 44463
$ws.Cells.Item(14,10).Value = 40
$ws.Cells.Item(14,11).Value = 2500
$ws.Cells.Item(14,12).Value = 2500
$ws.Cells.Item(14,13).Value = 2500
$ws.Cells.Item(14,16).Value = 2500

# Row 15
$ws.Cells.Item(15,4).Value = 44467
$ws.Cells.Item(15,10).Value = 50
$ws.Cells.Item(15,11).Value = 3000
$ws.Cells.Item(15,12).Value = 3000
$ws.Cells.Item(15,13).Value = 3000
$ws.Cells.Item(15,16).Value = 3000

# Row 16
$ws.Cells.Item(16,4).Value = 44473
$ws.Cells.Item(16,10).Value = 200
$ws.Cells.Item(16,11).Value = 1700
$ws.Cells.Item(16,12).Value = 1700
$ws.Cells.Item(16,13).Value = 1700
$ws.Cells.Item(16,16).Value = 1700

# Row 17
$ws.Cells.Item(17,4).Value = 44469
$ws.Cells.Item(17,10).Value = 1200
$ws.Cells.Item(17,11).Value = 1800
$ws.Cells.Item(17,12).Value = 1800
$ws.Cells.Item(17,13).Value = 1800
$ws.Cells.Item(17,16).Value = 1800

# Row 18
$ws.Cells.Item(18,4).Value = 44168
$ws.Cells.Item(18,10).Value = 150
$ws.Cells.Item(18,12).Value = 1000
$ws.Cells.Item(18,13).Value = 947
$ws.Cells.Item(18,16).Value = 947

# Row 19
$ws.Cells.Item(19,4).Value = 44161
$ws.Cells.Item(19,9).Value = "Primera"
$ws.Cells.Item(19,10).Value = 3000
$ws.Cells.Item(19,11).Value = 1000
$ws.Cells.Item(19,12).Value = 1000
$ws.Cells.Item(19,13).Value = 1000
$ws.Cells.Item(19,16).Value = 1000

# Row 20
$ws.Cells.Item(20,4).Value = 44165
$ws.Cells.Item(20,10).Value = 650
$ws.Cells.Item(20,11).Value = 900
$ws.Cells.Item(20,12).Value = 1100
$ws.Cells.Item(20,13).Value = 1008
$ws.Cells.Item(20,14).Value = "`$/kilo"
$ws.Cells.Item(20,16).Value = 1008
$ws.Cells.Item(20,17).Value = 1

# Row 21
$ws.Cells.Item(21,4).Value = 44165
$ws.Cells.Item(21,9).Value = "Segunda"
$ws.Cells.Item(21,10).Value = 180
$ws.Cells.Item(21,11).Value = 800
$ws.Cells.Item(21,12).Value = 800
$ws.Cells.Item(21,13).Value = 800
$ws.Cells.Item(21,16).Value = 800

# Row 22
$ws.Cells.Item(22,4).Value = 44475
$ws.Cells.Item(22,10).Value = 80
$ws.Cells.Item(22,11).Value = 17000
$ws.Cells.Item(22,12).Value = 17000
$ws.Cells.Item(22,13).Value = 17000
$ws.Cells.Item(22,14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(22,16).Value = 1700
$ws.Cells.Item(22,17).Value = 10

# Row 23
$ws.Cells.Item(23,4).Value = 44159
$ws.Cells.Item(23,10).Value = 2000
$ws.Cells.Item(23,11).Value = 1000
$ws.Cells.Item(23,12).Value = 1000
$ws.Cells.Item(23,13).Value = 1000
$ws.Cells.Item(23,16).Value = 1000

# Row 24
$ws.Cells.Item(24,4).Value = 44166
$ws.Cells.Item(24,9).Value = "Primera"
$ws.Cells.Item(24,10).Value = 285
$ws.Cells.Item(24,11).Value = 1000
$ws.Cells.Item(24,12).Value = 1100
$ws.Cells.Item(24,13).Value = 1054
$ws.Cells.Item(24,16).Value = 1054

# Row 25
$ws.Cells.Item(25,4).Value = 44476
$ws.Cells.Item(25,10).Value = 700
$ws.Cells.Item(25,11).Value = 1600
$ws.Cells.Item(25,12).Value = 1700
$ws.Cells.Item(25,13).Value = 1657
$ws.Cells.Item(25,16).Value = 1657

# Row 26
$ws.Cells.Item(26,4).Value = 44476
$ws.Cells.Item(26,10).Value = 100

# Row 27 (new)
$ws.Cells.Item(27,1).Value = 10
$ws.Cells.Item(27,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(27,3).Value = "La Araucanía"
$ws.Cells.Item(27,4).Value = 44466
$ws.Cells.Item(27,5).Value = 9
$ws.Cells.Item(27,6).Value = 300000000
$ws.Cells.Item(27,7).Value = "Espárragos"
$ws.Cells.Item(27,8).Value = "Sin especificar"
$ws.Cells.Item(27,9).Value = "Primera"
$ws.Cells.Item(27,10).Value = 300
$ws.Cells.Item(27,11).Value = 2000
$ws.Cells.Item(27,12).Value = 2000
$ws.Cells.Item(27,13).Value = 2000
$ws.Cells.Item(27,14).Value = "`$/kilo"
$ws.Cells.Item(27,15).Value = "Región del Maule"
$ws.Cells.Item(27,16).Value = 2000
$ws.Cells.Item(27,17).Value = 1
$ws.Cells.Item(27,18).Value = "Hortaliza"
$ws.Cells.Item(27,4).NumberFormat = $ws.Cells.Item(26,4).NumberFormat

# Row 28 (new)
$ws.Cells.Item(28,1).Value = 10
$ws.Cells.Item(28,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(28,3).Value = "La Araucanía"
$ws.Cells.Item(28,4).Value = 44466
$ws.Cells.Item(28,5).Value = 9
$ws.Cells.Item(28,6).Value = 300000000
$ws.Cells.Item(28,7).Value = "Espárragos"
$ws.Cells.Item(28,8).Value = "Sin especificar"
$ws.Cells.Item(28,9).Value = "Segunda"
$ws.Cells.Item(28,10).Value = 50
$ws.Cells.Item(28,11).Value = 1500
$ws.Cells.Item(28,12).Value = 1500
$ws.Cells.Item(28,13).Value = 1500
$ws.Cells.Item(28,14).Value = "`$/kilo"
$ws.Cells.Item(28,15).Value = "Región del Maule"
$ws.Cells.Item(28,16).Value = 1500
$ws.Cells.Item(28,17).Value = 1
$ws.Cells.Item(28,18).Value = "Hortaliza"
$ws.Cells.Item(28,4).NumberFormat = $ws.Cells.Item(27,4).NumberFormat

